# Natmi following Dr Hou advice:
# recompute the Cxcl16-Cxcr6 ligand-receptor table as a full 3x3
# cross-join of sending/target clusters (ECs, FAPs, sCs), expanding the
# sheet from 3 data rows to 9 data rows (rows 2-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl16"
$ws.Range("C2").Value = "Cxcr6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.703632666666667
$ws.Range("H2").Value = 11.110898
$ws.Range("I2").Value = 0.326750677445568
$ws.Range("J2").Value = 0.326750677445568
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.466623
$ws.Range("N2").Value = 4.399869000000001
$ws.Range("O2").Value = 0.4555124025226199
$ws.Range("P2").Value = 0.45551240252262
$ws.Range("Q2").Value = 5.431832852484668
$ws.Range("R2").Value = 48.88649567236201
$ws.Range("S2").Value = 0.1488389861091243
$ws.Range("T2").Value = 0.1488389861091243

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl16"
$ws.Range("C3").Value = "Cxcr6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.703632666666667
$ws.Range("H3").Value = 11.110898
$ws.Range("I3").Value = 0.326750677445568
$ws.Range("J3").Value = 0.326750677445568
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 1.071277666666667
$ws.Range("N3").Value = 3.213833
$ws.Range("O3").Value = 0.3327237222600216
$ws.Range("P3").Value = 0.3327237222600217
$ws.Range("Q3").Value = 3.967618961337112
$ws.Range("R3").Value = 35.70857065203401
$ws.Range("S3").Value = 0.1087177016506731
$ws.Range("T3").Value = 0.1087177016506731

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cxcl16"
$ws.Range("C4").Value = "Cxcr6"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.703632666666667
$ws.Range("H4").Value = 11.110898
$ws.Range("I4").Value = 0.326750677445568
$ws.Range("J4").Value = 0.326750677445568
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.6818206666666667
$ws.Range("N4").Value = 2.045462
$ws.Range("O4").Value = 0.2117638752173583
$ws.Range("P4").Value = 0.2117638752173583
$ws.Range("Q4").Value = 2.525213293875112
$ws.Range("R4").Value = 22.726919644876
$ws.Range("S4").Value = 0.06919398968577055
$ws.Range("T4").Value = 0.06919398968577056

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cxcl16"
$ws.Range("C5").Value = "Cxcr6"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7.123596666666667
$ws.Range("H5").Value = 21.37079
$ws.Range("I5").Value = 0.6284748640521197
$ws.Range("J5").Value = 0.6284748640521197
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.466623
$ws.Range("N5").Value = 4.399869000000001
$ws.Range("O5").Value = 0.4555124025226199
$ws.Range("P5").Value = 0.45551240252262
$ws.Range("Q5").Value = 10.44763071405667
$ws.Range("R5").Value = 94.02867642651002
$ws.Range("S5").Value = 0.286278095249458
$ws.Range("T5").Value = 0.2862780952494581

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cxcl16"
$ws.Range("C6").Value = "Cxcr6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 7.123596666666667
$ws.Range("H6").Value = 21.37079
$ws.Range("I6").Value = 0.6284748640521197
$ws.Range("J6").Value = 0.6284748640521197
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.071277666666667
$ws.Range("N6").Value = 3.213833
$ws.Range("O6").Value = 0.3327237222600216
$ws.Range("P6").Value = 0.3327237222600217
$ws.Range("Q6").Value = 7.631350015341112
$ws.Range("R6").Value = 68.68215013807
$ws.Range("S6").Value = 0.2091084961142823
$ws.Range("T6").Value = 0.2091084961142824

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cxcl16"
$ws.Range("C7").Value = "Cxcr6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 7.123596666666667
$ws.Range("H7").Value = 21.37079
$ws.Range("I7").Value = 0.6284748640521197
$ws.Range("J7").Value = 0.6284748640521197
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.6818206666666667
$ws.Range("N7").Value = 2.045462
$ws.Range("O7").Value = 0.2117638752173583
$ws.Range("P7").Value = 0.2117638752173583
$ws.Range("Q7").Value = 4.857015428331112
$ws.Range("R7").Value = 43.71313885498
$ws.Range("S7").Value = 0.1330882726883793
$ws.Range("T7").Value = 0.1330882726883793

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cxcl16"
$ws.Range("C8").Value = "Cxcr6"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5075066666666667
$ws.Range("H8").Value = 1.52252
$ws.Range("I8").Value = 0.04477445850231242
$ws.Range("J8").Value = 0.04477445850231242
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.466623
$ws.Range("N8").Value = 4.399869000000001
$ws.Range("O8").Value = 0.4555124025226199
$ws.Range("P8").Value = 0.45551240252262
$ws.Range("Q8").Value = 0.7443209499866668
$ws.Range("R8").Value = 6.698888549880001
$ws.Range("S8").Value = 0.02039532116403768
$ws.Range("T8").Value = 0.02039532116403768

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cxcl16"
$ws.Range("C9").Value = "Cxcr6"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5075066666666667
$ws.Range("H9").Value = 1.52252
$ws.Range("I9").Value = 0.04477445850231242
$ws.Range("J9").Value = 0.04477445850231242
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 1.071277666666667
$ws.Range("N9").Value = 3.213833
$ws.Range("O9").Value = 0.3327237222600216
$ws.Range("P9").Value = 0.3327237222600217
$ws.Range("Q9").Value = 0.5436805576844445
$ws.Range("R9").Value = 4.89312501916
$ws.Range("S9").Value = 0.01489752449506626
$ws.Range("T9").Value = 0.01489752449506626

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cxcl16"
$ws.Range("C10").Value = "Cxcr6"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5075066666666667
$ws.Range("H10").Value = 1.52252
$ws.Range("I10").Value = 0.04477445850231242
$ws.Range("J10").Value = 0.04477445850231242
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.6818206666666667
$ws.Range("N10").Value = 2.045462
$ws.Range("O10").Value = 0.2117638752173583
$ws.Range("P10").Value = 0.2117638752173583
$ws.Range("Q10").Value = 0.3460285338044445
$ws.Range("R10").Value = 3.114256804240001
$ws.Range("S10").Value = 0.009481612843208476
$ws.Range("T10").Value = 0.009481612843208477
